$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.283.54'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.11%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.930.51'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.01%  '

$ws.Range("E4").Value = '  +0.22%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '249.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.24%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7167'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.44%  '

$ws.Range("E7").Value = '  +0.24%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3215'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.80%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '27.37'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.46%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07110'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.49%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7921'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.37%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08033'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.19%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.927.07'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.18%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.372'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.79%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '94.73'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.01%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.60'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.80%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.273.61'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.09%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '257.20'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.25%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008068'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.66%  '

$ws.Range("E20").Value = '  -1.42%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.181.55'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.38%  '

$ws.Range("E22").Value = '  +0.21%  '

$ws.Range("E23").Value = '  +0.24%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.810'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.79%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.543'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.45%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.64'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.19'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.68%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.282'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.40%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1275'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.09%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.354'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.26%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.530'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.62%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.398'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.00%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.135'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.22%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05155'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.87%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.255'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.25%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7432'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.62%  '

$ws.Range("E37").Value = '  +0.57%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01954'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.74%  '

$ws.Range("E39").Value = '  -0.42%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '77.67'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.72%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.360'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.50%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4494'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.89%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.987'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.18%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8421'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.99%  '

$ws.Range("E45").Value = '  +0.18%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.02'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.78%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.716'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.26%  '

$ws.Range("E48").Value = '  +2.13%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.42'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.18%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06117'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.08%  '

$ws.Range("E51").Value = '  +2.51%  '
